$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("requirements (2)")
$ws2 = $wb.Worksheets.Item("requirements")
$ws3 = $wb.Worksheets.Item("officials")

# --- Fill in sample data on "requirements (2)" (A/O/X markers) ---
# Fix up the bottom row's interior border first so the cell style picks up
# the full thin border (matching the other data rows) before values go in.
$ws1.Range("C6:F6").Borders.Item(9).LineStyle = 1

$rows = 2,3,4,5,6
foreach ($r in $rows) {
  $ws1.Cells.Item($r, 2).Value = "A"
  $ws1.Cells.Item($r, 3).Value = "O"
  $ws1.Cells.Item($r, 4).Value = "O"
  $ws1.Cells.Item($r, 5).Value = "O"
  $ws1.Cells.Item($r, 6).Value = "O"
}
$ws1.Cells.Item(2, 11).Value = "X"
$ws1.Cells.Item(4, 11).Value = "X"
$ws1.Cells.Item(6, 11).Value = "X"

# --- New data validation: a dropdown list of A / O / X for the bulk of the
#     grid (rows 2-200), on top of the existing 0/1 validations used further
#     down the sheet ---
$ws1.Range("B2:AZ200").Validation.Add(3, 1, 1, """A, O, X""")

# --- Selections / active sheet ---
$null = $ws1.Range("K7").Select()
$null = $ws2.Range("C21").Select()
$null = $ws3.Range("E4:H4").Select()
